$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old per-row owner/member e-mail addresses with a single,
# shared "grjoh@jebosoft.onmicrosoft.com;tomjebo@jebosoft.onmicrosoft.com"
# string so the bot can match team/member mentions even with extra
# surrounding text (e.g. "<at>HuddleBot</at> ...").
$newVal = "grjoh@jebosoft.onmicrosoft.com;tomjebo@jebosoft.onmicrosoft.com"

$ws.Range("C2").Value = $newVal
$ws.Range("D2").Value = $newVal
$ws.Range("C3").Value = $newVal
$ws.Range("D3").Value = $newVal
$ws.Range("C4").Value = $newVal

# D4 was previously blank/untouched; give it the Hyperlink style so it
# matches the rest of the column formatting.
$ws.Range("D4").Style = "Hyperlink"

# Move / record the active selection like the author's session ended up.
$ws.Range("D13").Select()
